$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.919.58"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "1.636.90"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0637"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.864.06"
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").Value = "1.658.31"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.543"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.57%  "
$ws.Range("D16").Value = "0.0₃0764"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "25.949.49"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("E21").Value = "  -2.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.91"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.66%  "
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("E24").Value = "  +4.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0500"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("E34").Value = "  -4.40%  "
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.900"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("D37").Value = "1.132.29"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("E40").Value = "  -0.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.09%  "
$ws.Range("D44").Value = "1.773.52"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("E47").Value = "  +2.29%  "
$ws.Range("E48").Value = "  -0.50%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("E51").Value = "  -1.24%  "
